$wb = $excel.ActiveWorkbook

# --- Layer0 sheet ---
$ws0 = $wb.Worksheets.Item("Layer0")
$ws0.Range("B2").Value = 0.3948798872954241
$ws0.Range("C2").Value = -0.05516166515553721
$ws0.Range("B3").Value = -1.031384268746753
$ws0.Range("C3").Value = -0.5766246268025363
$ws0.Range("B4").Value = -1.558011315390599
$ws0.Range("C4").Value = -0.1213129715663307

# --- Layer1 sheet ---
$ws1 = $wb.Worksheets.Item("Layer1")
$ws1.Range("B2").Value = -0.3567884044806284
$ws1.Range("C2").Value = -0.1912283726950495
$ws1.Range("B3").Value = -2.10415017445991
$ws1.Range("C3").Value = -0.5856542588741644
$ws1.Range("B4").Value = -0.3767432051075367
$ws1.Range("C4").Value = 0.8894612098163481
